$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Fill in zero values for previously-empty "Points completed" (col B)
#     cells in the back half of Sprint 2 ---
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0

# End-of-sprint burst of completed points - this ripples through the
# cumulative "Points completed" (F) and "Points left" (G) shared formulas
# for the remaining rows (27-30).
$ws.Range("B27").Value = 4

# --- New "Target" marker in column H, styled like the section headers ---
$ws.Range("H20").Value = "Target"
$ws.Range("H20").Style = "Titolo 1"

$ws.Range("H21").Value = 21

# --- Selection now sits on H23 ---
$ws.Range("H23").Select()
